$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 1933.3334
$ws.Range("I34").Value = 1933.3334
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1933.3334
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1730.3334
$ws.Range("N34").ClearContents()
$ws.Range("H36").Value = 1933.3334
$ws.Range("I36").Value = 1933.3334
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 1933.3334
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -1218.3334
$ws.Range("N36").ClearContents()
$ws.Range("H112").Value = 1129.4
$ws.Range("I112").Value = 550
$ws.Range("J112").Value = 1149.3793
$ws.Range("K112").Value = 1650
$ws.Range("L112").Value = 3448.1379
$ws.Range("M112").Value = -542
$ws.Range("N112").Value = -5664.1379
$ws.Range("H138").Value = 2618.3057
$ws.Range("I138").Value = 2456.6924
$ws.Range("J138").Value = 2653.9153
$ws.Range("K138").Value = 7370.0772
$ws.Range("L138").Value = 7961.7459
$ws.Range("M138").Value = -2230.0772
$ws.Range("N138").Value = -18241.7459

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 29695.822
$ws.Range("I32").Value = 5285.415
$ws.Range("K32").Value = 5285.415
$ws.Range("M32").Value = -4998.415
$ws.Range("H68").Value = 49795
$ws.Range("J68").Value = 49795
$ws.Range("L68").Value = 49795
$ws.Range("N68").Value = -51417
$ws.Range("H71").Value = 49795
$ws.Range("J71").Value = 49795
$ws.Range("L71").Value = 149385
$ws.Range("N71").Value = -157497
$ws.Range("H122").Value = 2180.3
$ws.Range("I122").Value = 2163.6428
$ws.Range("K122").Value = 6490.928400000001
$ws.Range("M122").Value = -4040.928400000001
$ws.Range("H132").Value = 2040.2858
$ws.Range("I132").Value = 1454.7222
$ws.Range("K132").Value = 4364.1666
$ws.Range("M132").Value = -1834.1666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 144570.78
$ws.Range("I105").Value = 112942.22
$ws.Range("J105").Value = 201502.2
$ws.Range("K105").Value = 112942.22
$ws.Range("L105").Value = 201502.2
$ws.Range("M105").Value = -111195.22
$ws.Range("N105").Value = -204996.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 5176.857
$ws.Range("I132").Value = 5566.5557
$ws.Range("J132").Value = 4475.4
$ws.Range("K132").Value = 16699.6671
$ws.Range("L132").Value = 13426.2
$ws.Range("M132").Value = -14169.6671
$ws.Range("N132").Value = -18486.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 24875.652
$ws.Range("J37").Value = 24875.652
$ws.Range("L37").Value = 74626.95599999999
$ws.Range("N37").Value = -74850.95599999999
$ws.Range("H131").Value = 768.9400000000001
$ws.Range("I131").Value = 388.1579
$ws.Range("J131").Value = 858.2593000000001
$ws.Range("K131").Value = 1164.4737
$ws.Range("L131").Value = 2574.7779
$ws.Range("M131").Value = 3875.5263
$ws.Range("N131").Value = -12654.7779

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 500
$ws.Range("J4").Value = 500
$ws.Range("L4").Value = 500
$ws.Range("N4").Value = -724
$ws.Range("H63").Value = 17428.572
$ws.Range("J63").Value = 17428.572
$ws.Range("L63").Value = 17428.572
$ws.Range("N63").Value = -18800.572
$ws.Range("H66").Value = 17428.572
$ws.Range("J66").Value = 17428.572
$ws.Range("L66").Value = 52285.716
$ws.Range("N66").Value = -59149.716
$ws.Range("H122").Value = 800
$ws.Range("I122").Value = 800
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2400
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 50
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 4364.636
$ws.Range("I132").Value = 3499.8333
$ws.Range("J132").Value = 5402.4
$ws.Range("K132").Value = 10499.4999
$ws.Range("L132").Value = 16207.2
$ws.Range("M132").Value = -7969.499899999999
$ws.Range("N132").Value = -21267.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1322.409
$ws.Range("I22").Value = 2597.8
$ws.Range("J22").Value = 947.2941
$ws.Range("K22").Value = 2597.8
$ws.Range("L22").Value = 947.2941
$ws.Range("M22").Value = -2302.8
$ws.Range("N22").Value = -1537.2941
$ws.Range("H27").Value = 1322.409
$ws.Range("I27").Value = 2597.8
$ws.Range("J27").Value = 947.2941
$ws.Range("K27").Value = 2597.8
$ws.Range("L27").Value = 947.2941
$ws.Range("M27").Value = -2490.8
$ws.Range("N27").Value = -1161.2941
$ws.Range("H61").Value = 3252.889
$ws.Range("I61").Value = 2249.5
$ws.Range("J61").Value = 3539.5715
$ws.Range("K61").Value = 2249.5
$ws.Range("L61").Value = 3539.5715
$ws.Range("M61").Value = -2047.5
$ws.Range("N61").Value = -3943.5715
$ws.Range("H64").Value = 18333.334
$ws.Range("J64").Value = 18333.334
$ws.Range("L64").Value = 18333.334
$ws.Range("N64").Value = -18783.334
$ws.Range("H67").Value = 18333.334
$ws.Range("J67").Value = 18333.334
$ws.Range("L67").Value = 18333.334
$ws.Range("N67").Value = -19893.334
$ws.Range("H113").Value = 3252.889
$ws.Range("I113").Value = 2249.5
$ws.Range("J113").Value = 3539.5715
$ws.Range("K113").Value = 2249.5
$ws.Range("L113").Value = 3539.5715
$ws.Range("M113").Value = -79.5
$ws.Range("N113").Value = -7879.5715
$ws.Range("H132").Value = 3324.5186
$ws.Range("I132").Value = 3728.8948
$ws.Range("J132").Value = 2364.125
$ws.Range("K132").Value = 11186.6844
$ws.Range("L132").Value = 7092.375
$ws.Range("M132").Value = -8656.6844
$ws.Range("N132").Value = -12152.375
$ws.Range("H136").Value = 1195.4849
$ws.Range("I136").Value = 1049.862
$ws.Range("J136").Value = 2251.25
$ws.Range("K136").Value = 3149.586
$ws.Range("L136").Value = 6753.75
$ws.Range("M136").Value = -599.5860000000002
$ws.Range("N136").Value = -11853.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 39125
$ws.Range("J63").Value = 39125
$ws.Range("L63").Value = 39125
$ws.Range("N63").Value = -40373
$ws.Range("H66").Value = 39125
$ws.Range("J66").Value = 39125
$ws.Range("L66").Value = 117375
$ws.Range("N66").Value = -123615
$ws.Range("H113").Value = 677.64703
$ws.Range("I113").Value = 447.77777
$ws.Range("J113").Value = 936.25
$ws.Range("K113").Value = 1343.33331
$ws.Range("L113").Value = 2808.75
$ws.Range("M113").Value = 826.66669
$ws.Range("N113").Value = -7148.75
$ws.Range("H126").Value = 1458.5333
$ws.Range("I126").Value = 1647.8
$ws.Range("K126").Value = 4943.4
$ws.Range("M126").Value = -2473.4
